# Generate Report for Handback
#
# This script updates the localization-status workbook to reflect a
# completed handback:
#   * the "Ready for handoff" status text becomes
#     "Handed back: in sync with en-US"
#   * the zh-cn and de-de sheets gain "Latest Target File" (F) and
#     "Latest Handback File" (G) values (with hyperlinks) for both
#     data rows
#   * the "Latest Handback DateTime" (H) column is stamped with the
#     handback timestamp for each language

$wb = $excel.ActiveWorkbook

$mdUrl    = "https://github.com/OpenLocalizationTest/oltest/blob/a965481335f3658e335caf8acb6969d20bfd11bd/e2e/9993c348-c562-422b-8d38-0d8a9c505173.md"
$mdName   = "9993c348-c562-422b-8d38-0d8a9c505173.md"

$zhXlfUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d7c7d03529b3645dfec2a3897d47c720e6c4c275/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/9993c348-c562-422b-8d38-0d8a9c505173.9b05a357e2cce61ad44c1ed11ac0fa2e5af751b5.zh-cn.xlf"
$zhXlfName = "9993c348-c562-422b-8d38-0d8a9c505173.9b05a357e2cce61ad44c1ed11ac0fa2e5af751b5.zh-cn.xlf"

$deXlfUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4528dd617fc9df9722d9510a63b62aef87b5cc7f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/9993c348-c562-422b-8d38-0d8a9c505173.9b05a357e2cce61ad44c1ed11ac0fa2e5af751b5.de-de.xlf"
$deXlfName = "9993c348-c562-422b-8d38-0d8a9c505173.9b05a357e2cce61ad44c1ed11ac0fa2e5af751b5.de-de.xlf"

$handedBackStatus = "Handed back: in sync with en-US"

$zhHandbackTime = "2016-03-22 07:09:15"
$deHandbackTime = "2016-03-22 07:09:21"

$hyperlinkColor = 15570276  # RGB(0x64,0x95,0xED) packed BGR, matches the workbook's existing HyperLink font

function Set-HandbackRow($ws, $row, $targetUrl, $targetName, $handbackUrl, $handbackName, $handbackTime) {

    $fCell = $ws.Cells.Item($row, 6)   # F - Latest Target File
    $gCell = $ws.Cells.Item($row, 7)   # G - Latest Handback File
    $hCell = $ws.Cells.Item($row, 8)   # H - Latest Handback DateTime

    $fCell.Value = $targetName
    $ws.Hyperlinks.Add($fCell, $targetUrl, "", "", $targetName) | Out-Null
    $fCell.Font.Underline = $true
    $fCell.Font.Color = $hyperlinkColor

    $gCell.Value = $handbackName
    $ws.Hyperlinks.Add($gCell, $handbackUrl, "", "", $handbackName) | Out-Null
    $gCell.Font.Underline = $true
    $gCell.Font.Color = $hyperlinkColor

    $hCell.NumberFormat = "yyyy-mm-dd HH:mm:ss"
    $hCell.Value = $handbackTime
}

# --- Update the "Ready for handoff" status everywhere it appears -----------
$ws1 = $wb.Worksheets.Item(1)   # Overview
$ws1.Range("B2").Value = $handedBackStatus
$ws1.Range("C2").Value = $handedBackStatus
$ws1.Range("B3").Value = $handedBackStatus
$ws1.Range("C3").Value = $handedBackStatus

$ws2 = $wb.Worksheets.Item(2)   # zh-cn
$ws2.Range("C2").Value = $handedBackStatus
$ws2.Range("C3").Value = $handedBackStatus

$ws3 = $wb.Worksheets.Item(3)   # de-de
$ws3.Range("C2").Value = $handedBackStatus
$ws3.Range("C3").Value = $handedBackStatus

# --- zh-cn sheet: fill in Latest Target File / Latest Handback File --------
Set-HandbackRow $ws2 2 $mdUrl $mdName $zhXlfUrl $zhXlfName $zhHandbackTime
Set-HandbackRow $ws2 3 $mdUrl $mdName $zhXlfUrl $zhXlfName $zhHandbackTime

# --- de-de sheet: fill in Latest Target File / Latest Handback File --------
Set-HandbackRow $ws3 2 $mdUrl $mdName $deXlfUrl $deXlfName $deHandbackTime
Set-HandbackRow $ws3 3 $mdUrl $mdName $deXlfUrl $deXlfName $deHandbackTime

Write-Host "Handback report generated."
